$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(421, '2024-07-18 10:43:22', 'Пользователь User 2 написал сообщение 2; 33333')
    ,@(422, '2024-07-18 10:43:41', 'Пользователь User 2 вошёл в систему (cookie value)')
    ,@(423, '2024-07-18 10:47:13', 'Пользователь User 2 написал сообщение 2: 33333')
    ,@(424, '2024-07-18 11:09:37', 'Пользователь User 2 написал сообщение 2:23231')
    ,@(425, '2024-07-18 11:09:49', 'Пользователь User 2 вошёл в систему (cookie value)')
    ,@(426, '2024-07-18 11:10:00', 'Пользователь User 2 вышел из системы (logout)')
    ,@(427, '2024-07-18 11:10:27', 'Пользователь User 1 вошёл в систему (log in)')
    ,@(428, '2024-07-18 11:14:23', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(429, '2024-07-18 11:15:25', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(430, '2024-07-18 11:17:03', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(431, '2024-07-18 11:22:54', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(432, '2024-07-18 11:24:11', 'Пользователь User 1 вышел из системы (logout)')
    ,@(433, '2024-07-18 11:24:45', 'Пользователь User 1 вошёл в систему (log in)')
    ,@(434, '2024-07-18 11:25:33', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(435, '2024-07-18 11:25:36', 'Пользователь User 1 вышел из системы (logout)')
    ,@(436, '2024-07-18 11:26:06', 'Пользователь User 1 вошёл в систему (log in)')
    ,@(437, '2024-07-18 11:26:16', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(438, '2024-07-18 11:26:28', 'Пользователь User 1 написал сообщение 1:7777777')
    ,@(439, '2024-07-18 11:26:56', 'Пользователь User 1 вышел из системы (logout)')
    ,@(440, '2024-07-18 11:27:11', 'Пользователь User 2 вошёл в систему (log in)')
    ,@(441, '2024-07-18 11:28:19', 'Пользователь  вышел из системы (logout)')
    ,@(442, '2024-07-18 11:28:32', 'Пользователь User 1 вошёл в систему (log in)')
    ,@(443, '2024-07-18 11:28:39', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(444, '2024-07-18 11:28:58', 'Пользователь User 1 написал сообщение 1: вапвапвап')
    ,@(445, '2024-07-18 11:34:52', 'Пользователь User 1 вышел из системы (logout)')
    ,@(446, '2024-07-18 11:35:00', 'Пользователь User 2 вошёл в систему (log in)')
    ,@(447, '2024-07-18 11:35:07', 'Пользователь User 2 вышел из системы (logout)')
    ,@(448, '2024-07-18 11:35:33', 'Пользователь User 15 успешно прошёл регистрацию)')
    ,@(449, '2024-07-18 11:35:33', 'Пользователь User 15 вошёл в систему (cookie value)')
    ,@(450, '2024-07-18 11:35:50', 'Пользователь User 2 написал сообщение 15: ррркнуг')
    ,@(451, '2024-07-18 11:36:19', 'Пользователь User 15 вышел из системы (logout)')
    ,@(452, '2024-07-18 11:39:04', 'Пользователь User 16 успешно прошёл регистрацию)')
    ,@(453, '2024-07-18 11:39:16', 'Пользователь User 16 вошёл в систему (cookie value)')
    ,@(454, '2024-07-18 11:39:27', 'Пользователь User 2 написал сообщение 16: вапвпвпа')
    ,@(455, '2024-07-18 11:39:41', 'Пользователь User 16 вошёл в систему (cookie value)')
    ,@(456, '2024-07-18 11:39:46', 'Пользователь User 2 написал сообщение 16: кенкен')
    ,@(457, '2024-07-18 11:39:58', 'Пользователь User 16 вышел из системы (logout)')
    ,@(458, '2024-07-18 11:40:07', 'Пользователь User 16 вошёл в систему (log in)')
    ,@(459, '2024-07-18 11:40:10', 'Пользователь User 16 вошёл в систему (cookie value)')
    ,@(460, '2024-07-18 11:40:14', 'Пользователь User 16 написал сообщение 16 кенкен')
    ,@(461, '2024-07-18 11:40:25', 'Пользователь User 16 вышел из системы (logout)')
    ,@(462, '2024-07-18 11:40:31', 'Пользователь User 1 вошёл в систему (log in)')
    ,@(463, '2024-07-18 11:40:37', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(464, '2024-07-18 11:40:43', 'Пользователь User 1 написал сообщение 1 кеуке')
    ,@(465, '2024-07-18 11:40:49', 'Пользователь User 1 вышел из системы (logout)')
    ,@(466, '2024-07-18 11:40:58', 'Пользователь User 16 вошёл в систему (log in)')
    ,@(467, '2024-07-18 11:41:12', 'Пользователь User 16 вышел из системы (logout)')
    ,@(468, '2024-07-18 11:42:18', 'Пользователь User 2 вошёл в систему (log in)')
    ,@(469, '2024-07-18 11:42:26', 'Пользователь User 2 вошёл в систему (cookie value)')
    ,@(470, '2024-07-18 11:42:31', 'Пользователь User 2 написал сообщение 2 кн')
    ,@(471, '2024-07-18 11:42:44', 'Пользователь User 2 вышел из системы (logout)')
    ,@(472, '2024-07-18 11:42:55', 'Пользователь User 1 вошёл в систему (log in)')
    ,@(473, '2024-07-18 11:43:01', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(474, '2024-07-18 11:43:07', 'Пользователь User 1 написал сообщение 1 вапро')
    ,@(475, '2024-07-18 11:43:22', 'Пользователь User 1 вошёл в систему (cookie value)')
    ,@(476, '2024-07-18 11:43:27', 'Пользователь User 1 написал сообщение 1 ыкне')
    ,@(477, '2024-07-18 11:43:33', 'Пользователь User 1 вышел из системы (logout)')
    ,@(478, '2024-07-18 11:43:45', 'Пользователь User 17 успешно прошёл регистрацию)')
    ,@(479, '2024-07-18 11:44:01', 'Пользователь User 17 вошёл в систему (cookie value)')
    ,@(480, '2024-07-18 11:44:50', 'Пользователь User 17 вышел из системы (logout)')
    ,@(481, '2024-07-18 11:45:18', 'Пользователь  вышел из системы (logout)')
    ,@(482, '2024-07-18 11:45:29', 'Пользователь User 18 успешно прошёл регистрацию)')
    ,@(483, '2024-07-18 11:45:41', 'Пользователь User 18 вошёл в систему (log in)')
    ,@(484, '2024-07-18 11:45:45', 'Пользователь User 18 вошёл в систему (cookie value)')
    ,@(485, '2024-07-18 11:45:53', 'Пользователь User 18 написал сообщение 18: врповпрвпр')
    ,@(486, '2024-07-18 11:46:38', 'Пользователь User 18 вышел из системы (logout)')
)

foreach ($row in $data) {
    $r = $row[0]
    $dt = $row[1]
    $msg = $row[2]
    $ws.Cells.Item($r, 1).Value = $dt
    $ws.Cells.Item($r, 2).Value = $msg
}

Write-Output "Done adding rows"
